$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 5799.800999999999
$ws.Cells.Item(3, 2).Value = 5649.136500000001
$ws.Cells.Item(4, 2).Value = 5554.8325
$ws.Cells.Item(5, 2).Value = 5481.98
$ws.Cells.Item(6, 2).Value = 5456.8815
$ws.Cells.Item(7, 2).Value = 5512.395
$ws.Cells.Item(8, 2).Value = 5765.477227722772
$ws.Cells.Item(9, 2).Value = 6519.4465
$ws.Cells.Item(10, 2).Value = 7901.5895
$ws.Cells.Item(11, 2).Value = 13052.0005
$ws.Cells.Item(12, 2).Value = 14886.7985
$ws.Cells.Item(13, 2).Value = 14810.005
$ws.Cells.Item(14, 2).Value = 14699.3385
$ws.Cells.Item(15, 2).Value = 15357.636
$ws.Cells.Item(16, 2).Value = 15432.207
$ws.Cells.Item(17, 2).Value = 15338.9075
$ws.Cells.Item(18, 2).Value = 15792.80217625723
$ws.Cells.Item(19, 2).Value = 15557.50171551809
$ws.Cells.Item(20, 2).Value = 15152.49877462994
$ws.Cells.Item(21, 2).Value = 13812.90902852661
$ws.Cells.Item(22, 2).Value = 12007.84432898735
$ws.Cells.Item(23, 2).Value = 9523.143
$ws.Cells.Item(24, 2).Value = 6375.7855
$ws.Cells.Item(25, 2).Value = 5494.996500000001
